# Apply a cyclic rotation of species-observation data among rows 6, 7, 9 and 10
# (row 8 is unchanged). The row-specific columns (C, I, J, K, N, P, S:Y, Z, AA,
# AB, AD:AG, AT, AW:AY) stay put; only the data columns A, B, D, E, F, G, H, Q, R
# move: new row6 <- old row7, new row7 <- old row9, new row9 <- old row10,
# new row10 <- old row6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Capture current ("before") values for the affected rows/columns.
$rows = @(6, 7, 9, 10)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($c in $cols) {
        $snapshot[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Mapping describing which row's old data becomes the new data for each row.
$source = @{ 6 = 7; 7 = 9; 9 = 10; 10 = 6 }

foreach ($r in $rows) {
    $src = $source[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $snapshot[$src][$c]
    }
}

# Row 6 previously had an (empty) value in column L while rows 7, 9, 10 did not;
# after the rotation that empty cell moves from L6 to L10.
$ws.Range("L6").Value2 = $null
$ws.Range("L10").Value2 = ""
